$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column H, matching the formatting of the existing
# header cells (e.g. G1: bold, centered, bordered style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data values for the "Save" column.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
